# "The second commit with Delete data and save customer ID"
#
# 1. Update the existing customer email address on the "Customer" sheet
#    (the old email is replaced with a new one).
# 2. Create a new "CusID" worksheet (after "Customer") that stores a
#    Customer Name / Email Id / Phone / Customer ID table, saving the
#    previous customer's name, (old) email and phone number.

$wb = $excel.ActiveWorkbook

# --- Add the new "CusID" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cusIdSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$cusIdSheet.Name = "CusID"

# Header row
$cusIdSheet.Range("A1").Value = "Customer Name "
$cusIdSheet.Range("B1").Value = "Email Id "
$cusIdSheet.Range("C1").Value = "Phone"
$cusIdSheet.Range("D1").Value = "Customer ID"

# Saved customer record (name, old email, phone) - Customer ID left blank
$cusIdSheet.Range("A2").Value = "Matheq"
$cusIdSheet.Range("B2").Value = "rdete.rei654657@gmail.com"
$cusIdSheet.Range("C2").Value = 8174470105

# --- Update Customer sheet: replace the stored email address ---
$custSheet = $wb.Worksheets.Item("Customer")
$custSheet.Range("I2").Value = "rdesfgate657@gmail.com"

# Cosmetic touches matching the recorded column widths / view state
$cusIdSheet.Columns.Item(1).ColumnWidth = 14.166666666666666
$cusIdSheet.Columns.Item(2).ColumnWidth = 24.333333333333332
$cusIdSheet.Columns.Item(3).ColumnWidth = 10.333333333333334
$cusIdSheet.Columns.Item(4).ColumnWidth = 10.333333333333334

$cusIdSheet.Activate() | Out-Null
$cusIdSheet.Range("C2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 129
